$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph (paragraph 2, right after the title).
$metaPara = $d.Paragraphs.Item(2)
[void]$metaPara.Range.Delete()

# 2. Replace the final two-paragraph tail: the last paragraph currently holds the
#    "Create a feature image..." image-prompt text (italic). We need to:
#      - insert a new bold paragraph "Play Blazin Hot 7s Stack Em Up Slot for Free | Review 2021"
#        right before it
#      - change the last paragraph's text to the meta-description sentence, keeping italics
$count = $d.Paragraphs.Count
$secondLastP = $d.Paragraphs.Item($count - 1)
$lastP = $d.Paragraphs.Item($count)
$targetRange = $d.Range($secondLastP.Range.End, $lastP.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Blazin Hot 7s Stack Em Up Slot for Free | Review 2021</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our unbiased review of Blazin Hot 7s Stack Em Up slot. Learn how to play the game and try it for free. Discover pros and cons and RTP rate.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$targetRange.InsertXML($xml)
